# Applies the "Updated cryptos list" data refresh described by the commit diff.
# D-column "Price" cells are text (not numbers) in the original sheet, so we
# force text interpretation with a leading apostrophe (Excel's standard text-entry
# prefix) and then reset the resulting cell style back to "Normal" so no stray
# number-format / quote-prefix style is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell "D2" "57.743.53"
$ws.Range("E2").Value = "  +2.45%  "
Set-TextCell "D3" "3.060.55"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextCell "D5" "519.57"
$ws.Range("E5").Value = "  +2.84%  "
Set-TextCell "D6" "141.92"
$ws.Range("E6").Value = "  +3.64%  "
Set-TextCell "D7" "1.00"
$ws.Range("E7").Value = "  -0.03%  "
Set-TextCell "D8" "0.434"
$ws.Range("E8").Value = "  +1.40%  "
Set-TextCell "D9" "7.25"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +0.18%  "
Set-TextCell "D11" "0.376"
$ws.Range("E11").Value = "  +3.11%  "
Set-TextCell "D12" "3.589.08"
$ws.Range("E12").Value = "  +2.84%  "
$ws.Range("E13").Value = "  +3.33%  "
Set-TextCell "D14" "25.76"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("E15").Value = "  +0.32%  "
Set-TextCell "D16" "57.794.66"
$ws.Range("E16").Value = "  +2.61%  "
Set-TextCell "D17" "3.065.18"
$ws.Range("E17").Value = "  +2.91%  "
$ws.Range("E18").Value = "  +1.70%  "
Set-TextCell "D19" "12.84"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +1.20%  "
Set-TextCell "D21" "330.01"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  -0.01%  "
Set-TextCell "D23" "0.498"
$ws.Range("E23").Value = "  +1.29%  "
Set-TextCell "D24" "65.65"
$ws.Range("E24").Value = "  +1.69%  "
Set-TextCell "D25" "0.170"
$ws.Range("E25").Value = "  +3.97%  "
Set-TextCell "D26" "0.999"
$ws.Range("E26").Value = "  -0.07%  "
Set-TextCell "D27" "0.0₃0899"
$ws.Range("E27").Value = "  -1.89%  "
Set-TextCell "D28" "6.36"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  +4.07%  "
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D31" "1.19"
$ws.Range("E31").Value = "  +3.51%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D32" "20.66"
$ws.Range("E32").Value = "  +2.31%  "
Set-TextCell "D33" "154.96"
$ws.Range("E33").Value = "  +1.28%  "
Set-TextCell "D34" "4.50"
$ws.Range("E34").Value = "  +1.03%  "
Set-TextCell "D35" "27.27"
$ws.Range("E35").Value = "  +4.76%  "
Set-TextCell "D36" "5.95"
$ws.Range("E36").Value = "  +3.03%  "
$ws.Range("E37").Value = "  +0.80%  "
Set-TextCell "D38" "0.0673"
$ws.Range("E38").Value = "  +2.23%  "
Set-TextCell "D39" "3.106.23"
$ws.Range("E39").Value = "  +2.95%  "
Set-TextCell "D40" "3.92"
$ws.Range("E40").Value = "  +3.58%  "
Set-TextCell "D41" "36.69"
$ws.Range("E41").Value = "  -0.54%  "
Set-TextCell "D42" "1.00"
$ws.Range("E42").Value = "  -0.09%  "
Set-TextCell "D43" "0.652"
$ws.Range("E43").Value = "  +0.33%  "
Set-TextCell "D44" "2.267.44"
$ws.Range("E44").Value = "  +4.15%  "
Set-TextCell "D45" "0.0257"
$ws.Range("E45").Value = "  +10.17%  "
Set-TextCell "D46" "20.80"
$ws.Range("E46").Value = "  +6.85%  "
$ws.Range("E47").Value = "  +0.84%  "
Set-TextCell "D48" "5.87"
$ws.Range("E48").Value = "  +1.20%  "
Set-TextCell "D49" "0.921"
$ws.Range("E49").Value = "  +0.24%  "
Set-TextCell "D50" "261.33"
$ws.Range("E50").Value = "  +15.01%  "
Set-TextCell "D51" "0.724"
$ws.Range("E51").Value = "  +7.86%  "
